$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = 17
$ws.Range("H1").Value = 18
$ws.Range("I1").Value = 19

$ws.Range("I1").Select()
